$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ColorVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# --- Row 2: E2 and F2 get brand-new colors ---
# E2: #fefefe -> #fbfbfb  (font FF040404, fill FBFBFB)
$ws.Range("E2").Value = "#fbfbfb"
$ws.Range("E2").Font.Color = ColorVal 4 4 4
$ws.Range("E2").Interior.Color = ColorVal 251 251 251

# F2: #fafafa -> #ffffff  (font FF000000, fill FFFFFF)
$ws.Range("F2").Value = "#ffffff"
$ws.Range("F2").Font.Color = ColorVal 0 0 0
$ws.Range("F2").Interior.Color = ColorVal 255 255 255

# --- Row 3: E3 and F3 swap formatting/values ---
# Old E3 = "#fcd602" (font FF0329FD / fill FCD602) -> becomes F3's new value/format
# Old F3 = "#fefefe" (font FF010101 / fill FEFEFE) -> becomes E3's new value/format
# Copy E3's current (pre-edit) format onto F3 first, then copy E5's format (same as
# old F3: font FF010101 / fill FEFEFE) onto E3, preserving/reusing existing styles.
$ws.Range("E3").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null

$ws.Range("E5").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("E3").Value = "#fefefe"
$ws.Range("F3").Value = "#fcd602"
